$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A75").Value = "2025/12/05 21:00"
$ws.Range("B75").Value = "-"
$ws.Range("C75").Value = "-"
$ws.Range("D75").Value = "-"
$ws.Range("E75").Value = "-"
$ws.Range("F75").Value = "-"
$ws.Range("G75").Value = "-"
